$wb = $excel.ActiveWorkbook

# --- "dependencies" sheet: fix row9/row10 ETL dependency mapping ---
$deps = $wb.Worksheets.Item("dependencies")

# Row 9 (far_plots) columns K:O were wrongly set to "event-plot";
# they should reference the new "event-ext_id" parameter.
$deps.Range("K9:O9").Value = "event-ext_id"

# Row 10 (far_production_events) columns K:O were missing the
# "event-plot" relationship that had been mistakenly placed on row 9.
$deps.Range("K10:O10").Value = "event-plot"

# Update the active selection on the "dependencies" sheet.
$deps.Activate()
$deps.Range("K10").Select()

# --- "global" sheet: update the active selection ---
$glob = $wb.Worksheets.Item("global")
$glob.Activate()
$glob.Range("B6").Select()
